$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BEPEfCT")
$ws.Range("B2").Value = 1
